$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: for both locale sheets (zh-cn, de-de) the
# "Latest Target File" (F) / "Latest Handback File" (G) columns get filled in
# with the same file references already shown in "Source File Name" (A) and
# "Latest Handoff File" (D) for that row, the Status column (C) flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# "Latest Handback DateTime" (H) column gets a real timestamp.
# ---------------------------------------------------------------------------

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/306593eb8625b8c75cd18bcc4d553b800e32c3a5/e2e/e8b0186b-a8e4-4821-907b-ae5a903d2a7d.md"
$mdName   = "e8b0186b-a8e4-4821-907b-ae5a903d2a7d.md"

$statusText = "Handed back: in sync with en-US"

# --- zh-cn sheet (rows 2 and 3 both concern e8b0186b...md per the source
#     report, hence both reuse the same "Latest Target File" link) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b176567188d46cc45caf46748cfb99f95656a659/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e8b0186b-a8e4-4821-907b-ae5a903d2a7d.de989de1c6f941870db4db0d4dab7f11bb26bad7.zh-cn.xlf"
$zhXlfName   = "e8b0186b-a8e4-4821-907b-ae5a903d2a7d.de989de1c6f941870db4db0d4dab7f11bb26bad7.zh-cn.xlf"
$zhHandback  = "2016-03-25 01:27:46"

$wsZh.Cells.Item(2, 3).Value = $statusText
$wsZh.Cells.Item(2, 6).Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 6), $mdTarget, "", "", $mdName)
$wsZh.Cells.Item(2, 6).Style = "HyperLink"
$wsZh.Cells.Item(2, 7).Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 7), $zhXlfTarget, "", "", $zhXlfName)
$wsZh.Cells.Item(2, 7).Style = "HyperLink"
$wsZh.Cells.Item(2, 8).Value = $zhHandback

$wsZh.Cells.Item(3, 3).Value = $statusText
$wsZh.Cells.Item(3, 6).Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 6), $mdTarget, "", "", $mdName)
$wsZh.Cells.Item(3, 6).Style = "HyperLink"
$wsZh.Cells.Item(3, 7).Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 7), $zhXlfTarget, "", "", $zhXlfName)
$wsZh.Cells.Item(3, 7).Style = "HyperLink"
$wsZh.Cells.Item(3, 8).Value = $zhHandback

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94e3dcdcf084e6433c41f5865d27f89284ce489f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e8b0186b-a8e4-4821-907b-ae5a903d2a7d.de989de1c6f941870db4db0d4dab7f11bb26bad7.de-de.xlf"
$deXlfName   = "e8b0186b-a8e4-4821-907b-ae5a903d2a7d.de989de1c6f941870db4db0d4dab7f11bb26bad7.de-de.xlf"
$deHandback  = "2016-03-25 01:27:53"

$wsDe.Cells.Item(2, 3).Value = $statusText
$wsDe.Cells.Item(2, 6).Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 6), $mdTarget, "", "", $mdName)
$wsDe.Cells.Item(2, 6).Style = "HyperLink"
$wsDe.Cells.Item(2, 7).Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 7), $deXlfTarget, "", "", $deXlfName)
$wsDe.Cells.Item(2, 7).Style = "HyperLink"
$wsDe.Cells.Item(2, 8).Value = $deHandback

$wsDe.Cells.Item(3, 3).Value = $statusText
$wsDe.Cells.Item(3, 6).Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 6), $mdTarget, "", "", $mdName)
$wsDe.Cells.Item(3, 6).Style = "HyperLink"
$wsDe.Cells.Item(3, 7).Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 7), $deXlfTarget, "", "", $deXlfName)
$wsDe.Cells.Item(3, 7).Style = "HyperLink"
$wsDe.Cells.Item(3, 8).Value = $deHandback
